$d = $word.ActiveDocument

# Update the date header
$d.Content.Find.Execute("2025-05-13 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-14 Wednesday", 2) | Out-Null

$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "72+16=88"
$t.Cell(1, 2).Range.Text = "63-2=61"
$t.Cell(1, 3).Range.Text = "22+7=29"
$t.Cell(1, 4).Range.Text = "6+36=42"
$t.Cell(1, 5).Range.Text = "46-1=45"
$t.Cell(2, 1).Range.Text = "51+16=67"
$t.Cell(2, 2).Range.Text = "19+24=43"
$t.Cell(2, 3).Range.Text = "64-26=38"
$t.Cell(2, 4).Range.Text = "3+40=43"
$t.Cell(2, 5).Range.Text = "49+25=74"
$t.Cell(3, 1).Range.Text = "12+27=39"
$t.Cell(3, 2).Range.Text = "63-40=23"
$t.Cell(3, 3).Range.Text = "76-1=75"
$t.Cell(3, 4).Range.Text = "74-47=27"
$t.Cell(3, 5).Range.Text = "94-60=34"
$t.Cell(4, 1).Range.Text = "42+46=88"
$t.Cell(4, 2).Range.Text = "13-0=13"
$t.Cell(4, 3).Range.Text = "81-45=36"
$t.Cell(4, 4).Range.Text = "55-48=7"
$t.Cell(4, 5).Range.Text = "61-19=42"
$t.Cell(5, 1).Range.Text = "97-38=59"
$t.Cell(5, 2).Range.Text = "72-0=72"
$t.Cell(5, 3).Range.Text = "20+42=62"
$t.Cell(5, 4).Range.Text = "57+12=69"
$t.Cell(5, 5).Range.Text = "72-69=3"
$t.Cell(6, 1).Range.Text = "63-13=50"
$t.Cell(6, 2).Range.Text = "50-2=48"
$t.Cell(6, 3).Range.Text = "39+29=68"
$t.Cell(6, 4).Range.Text = "6+67=73"
$t.Cell(6, 5).Range.Text = "96-73=23"
$t.Cell(7, 1).Range.Text = "33-6=27"
$t.Cell(7, 2).Range.Text = "34+10=44"
$t.Cell(7, 3).Range.Text = "19-5=14"
$t.Cell(7, 4).Range.Text = "38+16=54"
$t.Cell(7, 5).Range.Text = "80-37=43"
$t.Cell(8, 1).Range.Text = "80+6=86"
$t.Cell(8, 2).Range.Text = "65-26=39"
$t.Cell(8, 3).Range.Text = "97-6=91"
$t.Cell(8, 4).Range.Text = "54-52=2"
$t.Cell(8, 5).Range.Text = "28+55=83"
$t.Cell(9, 1).Range.Text = "81-80=1"
$t.Cell(9, 2).Range.Text = "42+37=79"
$t.Cell(9, 3).Range.Text = "59+38=97"
$t.Cell(9, 4).Range.Text = "75-39=36"
$t.Cell(9, 5).Range.Text = "78-2=76"
$t.Cell(10, 1).Range.Text = "22+67=89"
$t.Cell(10, 2).Range.Text = "93-89=4"
$t.Cell(10, 3).Range.Text = "5+32=37"
$t.Cell(10, 4).Range.Text = "58-40=18"
$t.Cell(10, 5).Range.Text = "7+9=16"
$t.Cell(11, 1).Range.Text = "26+41=67"
$t.Cell(11, 2).Range.Text = "60-8=52"
$t.Cell(11, 3).Range.Text = "67-8=59"
$t.Cell(11, 4).Range.Text = "91-57=34"
$t.Cell(11, 5).Range.Text = "63-40=23"
$t.Cell(12, 1).Range.Text = "60+8=68"
$t.Cell(12, 2).Range.Text = "2+56=58"
$t.Cell(12, 3).Range.Text = "31+27=58"
$t.Cell(12, 4).Range.Text = "38+4=42"
$t.Cell(12, 5).Range.Text = "92-54=38"
$t.Cell(13, 1).Range.Text = "81-11=70"
$t.Cell(13, 2).Range.Text = "49+15=64"
$t.Cell(13, 3).Range.Text = "66+30=96"
$t.Cell(13, 4).Range.Text = "14+50=64"
$t.Cell(13, 5).Range.Text = "34-1=33"
$t.Cell(14, 1).Range.Text = "66-50=16"
$t.Cell(14, 2).Range.Text = "28+71=99"
$t.Cell(14, 3).Range.Text = "35-14=21"
$t.Cell(14, 4).Range.Text = "98-63=35"
$t.Cell(14, 5).Range.Text = "60-52=8"
$t.Cell(15, 1).Range.Text = "48-41=7"
$t.Cell(15, 2).Range.Text = "99-10=89"
$t.Cell(15, 3).Range.Text = "3+70=73"
$t.Cell(15, 4).Range.Text = "17+12=29"
$t.Cell(15, 5).Range.Text = "82+4=86"
$t.Cell(16, 1).Range.Text = "25+13=38"
$t.Cell(16, 2).Range.Text = "63+33=96"
$t.Cell(16, 3).Range.Text = "44+54=98"
$t.Cell(16, 4).Range.Text = "67+25=92"
$t.Cell(16, 5).Range.Text = "96-23=73"
$t.Cell(17, 1).Range.Text = "87-62=25"
$t.Cell(17, 2).Range.Text = "15+78=93"
$t.Cell(17, 3).Range.Text = "2+24=26"
$t.Cell(17, 4).Range.Text = "82-74=8"
$t.Cell(17, 5).Range.Text = "99-34=65"
$t.Cell(18, 1).Range.Text = "83-38=45"
$t.Cell(18, 2).Range.Text = "93+6=99"
$t.Cell(18, 3).Range.Text = "42+19=61"
$t.Cell(18, 4).Range.Text = "60+15=75"
$t.Cell(18, 5).Range.Text = "33-13=20"
$t.Cell(19, 1).Range.Text = "96-72=24"
$t.Cell(19, 2).Range.Text = "37+15=52"
$t.Cell(19, 3).Range.Text = "95-80=15"
$t.Cell(19, 4).Range.Text = "78-75=3"
$t.Cell(19, 5).Range.Text = "80+6=86"
$t.Cell(20, 1).Range.Text = "14+4=18"
$t.Cell(20, 2).Range.Text = "10+69=79"
$t.Cell(20, 3).Range.Text = "22-3=19"
$t.Cell(20, 4).Range.Text = "4+73=77"
$t.Cell(20, 5).Range.Text = "87-8=79"

Write-Output "done"
